# "add a new story" — DungeonGismo.xlsx
#
# Row 11 (dungeon story id 45000008) previously pointed at an orphaned,
# unused "没有找到穷奇" / "bossqiongqi2" story stub. This edit replaces it
# with a brand-new story: "rosemaryfield" — "连续3次进入【迷迭香田】"
# (repeat count bumped from 1 to 3 to match the new description).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the tag/id cell (column G) first, then the description (column C),
# so the shared-string table gains the two new entries in the same order
# the source workbook has them ("rosemaryfield" before the Chinese blurb).
$ws.Range("G11").Value = "rosemaryfield"
$ws.Range("C11").Value = "连续3次进入【迷迭香田】"

# FinishContinueCount for this row goes from 1 to 3.
$ws.Range("J11").Value = 3

# Leave the selection on the cell that was actually edited.
$ws.Range("C11").Select()
